$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-19: only the "most recent charge end time" (column D) changes;
# the underlying station/terminal/time-string values (A/B/C) are untouched.
$ws.Range("D2").Value2 = 45976.339004629626
$ws.Range("D3").Value2 = 45976.339004629626
$ws.Range("D4").Value2 = 45976.339004629626
$ws.Range("D5").Value2 = 45976.339004629626
$ws.Range("D6").Value2 = 45976.339004629626
$ws.Range("D7").Value2 = 45976.339004629626
$ws.Range("D8").Value2 = 45976.339004629626
$ws.Range("D9").Value2 = 45976.339004629626
$ws.Range("D10").Value2 = 45976.339004629626
$ws.Range("D11").Value2 = 45976.339004629626
$ws.Range("D12").Value2 = 45976.339004629626
$ws.Range("D13").Value2 = 45976.339004629626
$ws.Range("D14").Value2 = 45976.339004629626
$ws.Range("D15").Value2 = 45976.339004629626
$ws.Range("D16").Value2 = 45976.339004629626
$ws.Range("D17").Value2 = 45976.339004629626
$ws.Range("D18").Value2 = 45976.339004629626
$ws.Range("D19").Value2 = 45976.339004629626

# Rows 20-50: new/updated "not-charging" session records (A,B,C,D all set).
$ws.Range("A20").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value2 = "503号直流"
$ws.Range("C20").Value2 = 45973.5340625
$ws.Range("D20").Value2 = 45976.339004629626
$ws.Range("A21").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value2 = "802号直流"
$ws.Range("C21").Value2 = 45973.53917824074
$ws.Range("D21").Value2 = 45976.339004629626
$ws.Range("A22").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B22").Value2 = "104号直流"
$ws.Range("C22").Value2 = 45974.568333333336
$ws.Range("D22").Value2 = 45976.339004629626
$ws.Range("A23").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B23").Value2 = "001A号直流"
$ws.Range("C23").Value2 = 45975.04269675926
$ws.Range("D23").Value2 = 45976.339004629626
$ws.Range("A24").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B24").Value2 = "103号直流"
$ws.Range("C24").Value2 = 45975.1144212963
$ws.Range("D24").Value2 = 45976.339004629626
$ws.Range("A25").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value2 = "704号直流"
$ws.Range("C25").Value2 = 45975.116527777776
$ws.Range("D25").Value2 = 45976.339004629626
$ws.Range("A26").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B26").Value2 = "401号直流"
$ws.Range("C26").Value2 = 45975.22027777778
$ws.Range("D26").Value2 = 45976.339004629626
$ws.Range("A27").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value2 = "101号直流"
$ws.Range("C27").Value2 = 45975.291817129626
$ws.Range("D27").Value2 = 45976.339004629626
$ws.Range("A28").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B28").Value2 = "203号直流"
$ws.Range("C28").Value2 = 45975.34092592593
$ws.Range("D28").Value2 = 45976.339004629626
$ws.Range("A29").Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B29").Value2 = "110号直流"
$ws.Range("C29").Value2 = 45975.42833333334
$ws.Range("D29").Value2 = 45976.339004629626
$ws.Range("A30").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B30").Value2 = "905号直流"
$ws.Range("C30").Value2 = 45975.506875
$ws.Range("D30").Value2 = 45976.339004629626
$ws.Range("A31").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B31").Value2 = "401号直流"
$ws.Range("C31").Value2 = 45975.53476851852
$ws.Range("D31").Value2 = 45976.339004629626
$ws.Range("A32").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B32").Value2 = "102号直流"
$ws.Range("C32").Value2 = 45975.53528935185
$ws.Range("D32").Value2 = 45976.339004629626
$ws.Range("A33").Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B33").Value2 = "205号直流"
$ws.Range("C33").Value2 = 45975.537824074076
$ws.Range("D33").Value2 = 45976.339004629626
$ws.Range("A34").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B34").Value2 = "705号直流"
$ws.Range("C34").Value2 = 45975.543807870374
$ws.Range("D34").Value2 = 45976.339004629626
$ws.Range("A35").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B35").Value2 = "101号直流"
$ws.Range("C35").Value2 = 45975.5471412037
$ws.Range("D35").Value2 = 45976.339004629626
$ws.Range("A36").Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B36").Value2 = "108号直流"
$ws.Range("C36").Value2 = 45975.55060185185
$ws.Range("D36").Value2 = 45976.339004629626
$ws.Range("A37").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B37").Value2 = "003B号直流"
$ws.Range("C37").Value2 = 45975.55443287037
$ws.Range("D37").Value2 = 45976.339004629626
$ws.Range("A38").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B38").Value2 = "904号直流"
$ws.Range("C38").Value2 = 45975.56085648148
$ws.Range("D38").Value2 = 45976.339004629626
$ws.Range("A39").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B39").Value2 = "702号直流"
$ws.Range("C39").Value2 = 45975.570185185185
$ws.Range("D39").Value2 = 45976.339004629626
$ws.Range("A40").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B40").Value2 = "905号直流"
$ws.Range("C40").Value2 = 45975.58256944444
$ws.Range("D40").Value2 = 45976.339004629626
$ws.Range("A41").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value2 = "404号直流"
$ws.Range("C41").Value2 = 45975.58453703704
$ws.Range("D41").Value2 = 45976.339004629626
$ws.Range("A42").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value2 = "804号直流"
$ws.Range("C42").Value2 = 45975.589641203704
$ws.Range("D42").Value2 = 45976.339004629626
$ws.Range("A43").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B43").Value2 = "502号直流"
$ws.Range("C43").Value2 = 45975.6040625
$ws.Range("D43").Value2 = 45976.339004629626
$ws.Range("A44").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B44").Value2 = "002B号直流"
$ws.Range("C44").Value2 = 45975.60476851852
$ws.Range("D44").Value2 = 45976.339004629626
$ws.Range("A45").Value2 = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B45").Value2 = "201号直流"
$ws.Range("C45").Value2 = 45975.66667824074
$ws.Range("D45").Value2 = 45976.339004629626
$ws.Range("A46").Value2 = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B46").Value2 = "009A号直流"
$ws.Range("C46").Value2 = 45975.69802083333
$ws.Range("D46").Value2 = 45976.339004629626
$ws.Range("A47").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B47").Value2 = "A02号直流"
$ws.Range("C47").Value2 = 45975.74300925926
$ws.Range("D47").Value2 = 45976.339004629626
$ws.Range("A48").Value2 = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B48").Value2 = "207号直流"
$ws.Range("C48").Value2 = 45975.74886574074
$ws.Range("D48").Value2 = 45976.339004629626
$ws.Range("A49").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B49").Value2 = "604号直流"
$ws.Range("C49").Value2 = 45975.764236111114
$ws.Range("D49").Value2 = 45976.339004629626
$ws.Range("A50").Value2 = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B50").Value2 = "B04号直流"
$ws.Range("C50").Value2 = 45975.76695601852
$ws.Range("D50").Value2 = 45976.339004629626

# Move the active selection to match the saved view state (F8).
$ws.Range("F8").Select()
